$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "332.52"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "1.04%"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "45.66"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "2.96%"

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "5.514"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "0.18%"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "0.08525"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "5.58%"

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "2.058"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "0.04%"

$c = $ws.Cells.Item(7, 2)
$c.NumberFormat = "@"
$c.Value = "GateToken"
$c = $ws.Cells.Item(7, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "4.443"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "0.81%"

$c = $ws.Cells.Item(8, 2)
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c = $ws.Cells.Item(8, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.9907"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "3.89%"

$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = "@"
$c.Value = "BTSEToken"
$c = $ws.Cells.Item(9, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "2.573"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "-0.51%"

$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = "@"
$c.Value = "LiechtensteinCryptoassetsExchange"
$c = $ws.Cells.Item(10, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.1154"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "1.44%"

$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = "@"
$c.Value = "WazirX"
$c = $ws.Cells.Item(11, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.1919"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "0.97%"

$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = "@"
$c.Value = "MCDex"
$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "9.465"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "-6.91%"

$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = "@"
$c.Value = "MandalaExchangeToken"
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.09755"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "-1.36%"

$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = "@"
$c.Value = "BitrueCoin"
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.04720"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "-3.31%"

$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = "@"
$c.Value = "BitMartToken"
$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.1058"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "-0.52%"

$c = $ws.Cells.Item(16, 2)
$c.NumberFormat = "@"
$c.Value = "BitForexToken"
$c = $ws.Cells.Item(16, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.001304"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "3.11%"

$c = $ws.Cells.Item(17, 2)
$c.NumberFormat = "@"
$c.Value = "TigerCash"
$c = $ws.Cells.Item(17, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.005921"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "1.72%"

$c = $ws.Cells.Item(18, 2)
$c.NumberFormat = "@"
$c.Value = "LEO"
$c = $ws.Cells.Item(18, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.383"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "0.17%"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.3356"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "-1.49%"

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "0.1374"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "-0.62%"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.2552"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "-1.00%"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.04147"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "1.45%"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "0.001301"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "-0.01%"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "0.004601"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "5.71%"

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "0.0001302"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "4.13%"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.0002987"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "-20.16%"

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.02760"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "6.50%"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.05719"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "-0.27%"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.007873"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "3.84%"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.1434"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "2.31%"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.007245"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "-1.55%"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.002155"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "7.29%"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.008101"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "-10.55%"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.3557"

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.00007057"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "0.61%"

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "0.13%"

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "0.27%"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.003456"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "-1.34%"

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "1.04%"

$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "0.13%"

$wb.Save()